# Auto-generated Excel COM-interop script applying numeric corrections
# to the Exodus_Profits.xlsx sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each sheet is an Excel Table of FFXIV leve crafting profit data; only specific
# numeric cells (H/I/J/K/L/M/N columns) are updated to refreshed market-board values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 2675.625  # H135
$ws.Cells.Item(135, 9).Value = 2675.625  # I135
$ws.Cells.Item(135, 11).Value = 24080.625  # K135
$ws.Cells.Item(135, 13).Value = -21545.625  # M135
$ws.Cells.Item(137, 8).Value = 560701.0600000001  # H137
$ws.Cells.Item(137, 9).Value = 2403.5  # I137
$ws.Cells.Item(137, 10).Value = 1039241.8  # J137
$ws.Cells.Item(137, 11).Value = 7210.5  # K137
$ws.Cells.Item(137, 12).Value = 3117725.4  # L137
$ws.Cells.Item(137, 13).Value = -4660.5  # M137
$ws.Cells.Item(137, 14).Value = -3122825.4  # N137
$ws.Cells.Item(141, 8).Value = 1936.5483  # H141
$ws.Cells.Item(141, 9).Value = 1587.9259  # I141
$ws.Cells.Item(141, 10).Value = 4289.75  # J141
$ws.Cells.Item(141, 11).Value = 4763.7777  # K141
$ws.Cells.Item(141, 12).Value = 12869.25  # L141
$ws.Cells.Item(141, 13).Value = 416.2223000000004  # M141
$ws.Cells.Item(141, 14).Value = -23229.25  # N141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(8, 8).Value = 99  # H8
$ws.Cells.Item(8, 10).Value = 99  # J8
$ws.Cells.Item(8, 12).Value = 99  # L8
$ws.Cells.Item(8, 14).Value = -387  # N8
$ws.Cells.Item(32, 8).Value = 8131.939  # H32
$ws.Cells.Item(32, 9).Value = 2949.8  # I32
$ws.Cells.Item(32, 10).Value = 16314.263  # J32
$ws.Cells.Item(32, 11).Value = 2949.8  # K32
$ws.Cells.Item(32, 12).Value = 16314.263  # L32
$ws.Cells.Item(32, 13).Value = -2662.8  # M32
$ws.Cells.Item(32, 14).Value = -16888.263  # N32
$ws.Cells.Item(44, 8).Value = 104144  # H44
$ws.Cells.Item(44, 9).Value = 86000  # I44
$ws.Cells.Item(44, 11).Value = 86000  # K44
$ws.Cells.Item(44, 13).Value = -85512  # M44
$ws.Cells.Item(61, 8).Value = 69141  # H61
$ws.Cells.Item(61, 9).Value = 2693.8462  # I61
$ws.Cells.Item(61, 11).Value = 2693.8462  # K61
$ws.Cells.Item(61, 13).Value = -2481.8462  # M61
$ws.Cells.Item(123, 8).Value = 52598  # H123
$ws.Cells.Item(123, 10).Value = 52598  # J123
$ws.Cells.Item(123, 12).Value = 52598  # L123
$ws.Cells.Item(123, 14).Value = -62398  # N123
$ws.Cells.Item(132, 8).Value = 2679.0476  # H132
$ws.Cells.Item(132, 9).Value = 2677.4666  # I132
$ws.Cells.Item(132, 11).Value = 8032.399800000001  # K132
$ws.Cells.Item(132, 13).Value = -5502.399800000001  # M132
$ws.Cells.Item(136, 8).Value = 69141  # H136
$ws.Cells.Item(136, 9).Value = 2693.8462  # I136
$ws.Cells.Item(136, 11).Value = 8081.5386  # K136
$ws.Cells.Item(136, 13).Value = -5531.5386  # M136
$ws.Cells.Item(138, 8).Value = 27499.5  # H138
$ws.Cells.Item(138, 10).Value = 24999  # J138
$ws.Cells.Item(138, 12).Value = 24999  # L138
$ws.Cells.Item(138, 14).Value = -35279  # N138

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 131723.08  # H20
$ws.Cells.Item(20, 9).Value = 163708.95  # I20
$ws.Cells.Item(20, 10).Value = 3779.6  # J20
$ws.Cells.Item(20, 11).Value = 163708.95  # K20
$ws.Cells.Item(20, 12).Value = 3779.6  # L20
$ws.Cells.Item(20, 13).Value = -163461.95  # M20
$ws.Cells.Item(20, 14).Value = -4273.6  # N20
$ws.Cells.Item(75, 8).Value = 11807.143  # H75
$ws.Cells.Item(75, 9).Value = 9358.333000000001  # I75
$ws.Cells.Item(75, 10).Value = 26500  # J75
$ws.Cells.Item(75, 11).Value = 9358.333000000001  # K75
$ws.Cells.Item(75, 12).Value = 26500  # L75
$ws.Cells.Item(75, 13).Value = -8422.333000000001  # M75
$ws.Cells.Item(75, 14).Value = -28372  # N75
$ws.Cells.Item(78, 8).Value = 11807.143  # H78
$ws.Cells.Item(78, 9).Value = 9358.333000000001  # I78
$ws.Cells.Item(78, 10).Value = 26500  # J78
$ws.Cells.Item(78, 11).Value = 28074.999  # K78
$ws.Cells.Item(78, 12).Value = 79500  # L78
$ws.Cells.Item(78, 13).Value = -23394.999  # M78
$ws.Cells.Item(78, 14).Value = -88860  # N78
$ws.Cells.Item(99, 8).Value = 1344614.8  # H99
$ws.Cells.Item(99, 9).Value = 50663.9  # I99
$ws.Cells.Item(99, 11).Value = 50663.9  # K99
$ws.Cells.Item(99, 13).Value = -49165.9  # M99
$ws.Cells.Item(126, 8).Value = 54287.5  # H126
$ws.Cells.Item(126, 10).Value = 54287.5  # J126
$ws.Cells.Item(126, 12).Value = 54287.5  # L126
$ws.Cells.Item(126, 14).Value = -64167.5  # N126
$ws.Cells.Item(134, 8).Value = 4506.037  # H134
$ws.Cells.Item(134, 9).Value = 2579.238  # I134
$ws.Cells.Item(134, 11).Value = 7737.714  # K134
$ws.Cells.Item(134, 13).Value = -5202.714  # M134
$ws.Cells.Item(138, 8).Value = 81547.44500000001  # H138
$ws.Cells.Item(138, 10).Value = 81547.44500000001  # J138
$ws.Cells.Item(138, 12).Value = 81547.44500000001  # L138
$ws.Cells.Item(138, 14).Value = -91827.44500000001  # N138

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1580.0476  # H58
$ws.Cells.Item(58, 9).Value = 1407  # I58
$ws.Cells.Item(58, 11).Value = 1407  # K58
$ws.Cells.Item(58, 13).Value = -1204  # M58
$ws.Cells.Item(94, 8).Value = 998.3333  # H94
$ws.Cells.Item(94, 9).Value = 1000  # I94
$ws.Cells.Item(94, 10).Value = 997.5  # J94
$ws.Cells.Item(94, 11).Value = 1000  # K94
$ws.Cells.Item(94, 12).Value = 997.5  # L94
$ws.Cells.Item(94, 13).Value = -549  # M94
$ws.Cells.Item(94, 14).Value = -1899.5  # N94
$ws.Cells.Item(136, 8).Value = 1580.0476  # H136
$ws.Cells.Item(136, 9).Value = 1407  # I136
$ws.Cells.Item(136, 11).Value = 4221  # K136
$ws.Cells.Item(136, 13).Value = -1671  # M136
$ws.Cells.Item(138, 8).Value = 67468  # H138
$ws.Cells.Item(138, 10).Value = 69408.89  # J138
$ws.Cells.Item(138, 12).Value = 69408.89  # L138
$ws.Cells.Item(138, 14).Value = -79688.89  # N138

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 4280.4  # H7
$ws.Cells.Item(7, 10).Value = 10375.5  # J7
$ws.Cells.Item(7, 12).Value = 31126.5  # L7
$ws.Cells.Item(7, 14).Value = -31350.5  # N7
$ws.Cells.Item(14, 8).Value = 84691.25  # H14
$ws.Cells.Item(14, 9).Value = 84691.25  # I14
$ws.Cells.Item(14, 11).Value = 254073.75  # K14
$ws.Cells.Item(14, 13).Value = -253900.75  # M14
$ws.Cells.Item(70, 8).Value = 0  # H70
$ws.Cells.Item(70, 9).Value = 0  # I70
$ws.Cells.Item(70, 11).Value = 0  # K70
$ws.Cells.Item(70, 13).ClearContents()  # M70
$ws.Cells.Item(73, 8).Value = 0  # H73
$ws.Cells.Item(73, 9).Value = 0  # I73
$ws.Cells.Item(73, 11).Value = 0  # K73
$ws.Cells.Item(73, 13).ClearContents()  # M73
$ws.Cells.Item(129, 8).Value = 41667160  # H129
$ws.Cells.Item(129, 9).Value = 560.2857  # I129
$ws.Cells.Item(129, 11).Value = 1680.8571  # K129
$ws.Cells.Item(129, 13).Value = 3319.1429  # M129
$ws.Cells.Item(131, 8).Value = 84460.25  # H131
$ws.Cells.Item(131, 9).Value = 111823.78  # I131
$ws.Cells.Item(131, 10).Value = 2369.6667  # J131
$ws.Cells.Item(131, 11).Value = 335471.34  # K131
$ws.Cells.Item(131, 12).Value = 7109.000100000001  # L131
$ws.Cells.Item(131, 13).Value = -330431.34  # M131
$ws.Cells.Item(131, 14).Value = -17189.0001  # N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 1628.4286  # H80
$ws.Cells.Item(80, 9).Value = 999.5  # I80
$ws.Cells.Item(80, 10).Value = 1880  # J80
$ws.Cells.Item(80, 11).Value = 999.5  # K80
$ws.Cells.Item(80, 12).Value = 1880  # L80
$ws.Cells.Item(80, 13).Value = -1.5  # M80
$ws.Cells.Item(80, 14).Value = -3876  # N80
$ws.Cells.Item(83, 8).Value = 1628.4286  # H83
$ws.Cells.Item(83, 9).Value = 999.5  # I83
$ws.Cells.Item(83, 10).Value = 1880  # J83
$ws.Cells.Item(83, 11).Value = 4997.5  # K83
$ws.Cells.Item(83, 12).Value = 9400  # L83
$ws.Cells.Item(83, 13).Value = -5.5  # M83
$ws.Cells.Item(83, 14).Value = -19384  # N83

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value = 1000  # H3
$ws.Cells.Item(3, 10).Value = 1000  # J3
$ws.Cells.Item(3, 12).Value = 1000  # L3
$ws.Cells.Item(3, 14).Value = -1224  # N3
$ws.Cells.Item(11, 8).Value = 1000  # H11
$ws.Cells.Item(11, 10).Value = 1000  # J11
$ws.Cells.Item(11, 12).Value = 1000  # L11
$ws.Cells.Item(11, 14).Value = -1280  # N11
$ws.Cells.Item(15, 8).Value = 1000  # H15
$ws.Cells.Item(15, 10).Value = 1000  # J15
$ws.Cells.Item(15, 12).Value = 1000  # L15
$ws.Cells.Item(15, 14).Value = -1340  # N15
$ws.Cells.Item(25, 8).Value = 1625  # H25
$ws.Cells.Item(25, 9).Value = 1500  # I25
$ws.Cells.Item(25, 10).Value = 1750  # J25
$ws.Cells.Item(25, 11).Value = 1500  # K25
$ws.Cells.Item(25, 12).Value = 1750  # L25
$ws.Cells.Item(25, 13).Value = -1270  # M25
$ws.Cells.Item(25, 14).Value = -2210  # N25
$ws.Cells.Item(55, 8).Value = 4878690.5  # H55
$ws.Cells.Item(55, 9).Value = 493.4  # I55
$ws.Cells.Item(55, 10).Value = 7693035  # J55
$ws.Cells.Item(55, 11).Value = 493.4  # K55
$ws.Cells.Item(55, 12).Value = 7693035  # L55
$ws.Cells.Item(55, 13).Value = -320.4  # M55
$ws.Cells.Item(55, 14).Value = -7693381  # N55
$ws.Cells.Item(136, 8).Value = 4156.7617  # H136
$ws.Cells.Item(136, 9).Value = 4610.857  # I136
$ws.Cells.Item(136, 11).Value = 13832.571  # K136
$ws.Cells.Item(136, 13).Value = -11282.571  # M136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 2500  # H18
$ws.Cells.Item(18, 10).Value = 2500  # J18
$ws.Cells.Item(18, 12).Value = 2500  # L18
$ws.Cells.Item(18, 14).Value = -2846  # N18
$ws.Cells.Item(80, 8).Value = 75000  # H80
$ws.Cells.Item(80, 10).Value = 75000  # J80
$ws.Cells.Item(80, 12).Value = 75000  # L80
$ws.Cells.Item(80, 14).Value = -76996  # N80
$ws.Cells.Item(83, 8).Value = 75000  # H83
$ws.Cells.Item(83, 10).Value = 75000  # J83
$ws.Cells.Item(83, 12).Value = 225000  # L83
$ws.Cells.Item(83, 14).Value = -234984  # N83
$ws.Cells.Item(95, 8).Value = 179497  # H95
$ws.Cells.Item(95, 10).Value = 179497  # J95
$ws.Cells.Item(95, 12).Value = 179497  # L95
$ws.Cells.Item(95, 14).Value = -184989  # N95
$ws.Cells.Item(100, 8).Value = 3247442.2  # H100
$ws.Cells.Item(100, 9).Value = 3968930.5  # I100
$ws.Cells.Item(100, 11).Value = 7937861  # K100
$ws.Cells.Item(100, 13).Value = -7937320  # M100
$ws.Cells.Item(108, 8).Value = 70284.5  # H108
$ws.Cells.Item(108, 9).Value = 40569  # I108
$ws.Cells.Item(108, 11).Value = 40569  # K108
$ws.Cells.Item(108, 13).Value = -36729  # M108
$ws.Cells.Item(123, 8).Value = 73000  # H123
$ws.Cells.Item(123, 10).Value = 73000  # J123
$ws.Cells.Item(123, 12).Value = 73000  # L123
$ws.Cells.Item(123, 14).Value = -82800  # N123
$ws.Cells.Item(136, 8).Value = 1843.9166  # H136
$ws.Cells.Item(136, 9).Value = 1430.7778  # I136
$ws.Cells.Item(136, 11).Value = 4292.3334  # K136
$ws.Cells.Item(136, 13).Value = -1742.3334  # M136

